$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text interpretation (avoids Excel
# auto-converting numeric-looking strings like "1.002" into numbers), and
# then restore the cell style back to Normal so no stray formatting is left
# behind (matches the original file which has no explicit style on these cells).
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "30.106.87"

# Row 3
Set-TextValue $ws.Range("D3") "1.926.52"
$ws.Range("E3").Value = "  +3.05%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.002"
$ws.Range("E4").Value = "  -0.60%  "

# Row 5
Set-TextValue $ws.Range("D5") "321.58"
$ws.Range("E5").Value = "  +2.03%  "

# Row 6
Set-TextValue $ws.Range("D6") "1.001"
$ws.Range("E6").Value = "  -0.59%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.5171"
$ws.Range("E7").Value = "  +2.32%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.3995"
$ws.Range("E8").Value = "  +2.66%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.08485"
$ws.Range("E9").Value = "  +1.69%  "

# Row 10
Set-TextValue $ws.Range("D10") "43.03"
$ws.Range("E10").Value = "  +3.07%  "

# Row 11
$ws.Range("E11").Value = "  +2.27%  "

# Row 12
Set-TextValue $ws.Range("B12") "Solana"
Set-TextValue $ws.Range("C12") "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue $ws.Range("D12") "21.29"
$ws.Range("E12").Value = "  +4.62%  "

# Row 13
Set-TextValue $ws.Range("B13") "Polkadot"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D13") "6.336"
$ws.Range("E13").Value = "  +2.30%  "

# Row 14
Set-TextValue $ws.Range("D14") "1.923.64"
$ws.Range("E14").Value = "  +2.83%  "

# Row 15
Set-TextValue $ws.Range("D15") "7.386"
$ws.Range("E15").Value = "  +2.23%  "

# Row 16
$ws.Range("E16").Value = "  -0.58%  "

# Row 17
Set-TextValue $ws.Range("D17") "94.49"
$ws.Range("E17").Value = "  +4.00%  "

# Row 18
Set-TextValue $ws.Range("D18") "0.00001117"
$ws.Range("E18").Value = "  +1.65%  "

# Row 19
Set-TextValue $ws.Range("D19") "0.06760"
$ws.Range("E19").Value = "  +1.07%  "

# Row 20
Set-TextValue $ws.Range("D20") "18.00"
$ws.Range("E20").Value = "  +2.10%  "

# Row 21
Set-TextValue $ws.Range("D21") "1.001"
$ws.Range("E21").Value = "  -0.62%  "

# Row 22
Set-TextValue $ws.Range("D22") "6.077"
$ws.Range("E22").Value = "  +2.94%  "

# Row 23
Set-TextValue $ws.Range("D23") "30.108.24"
$ws.Range("E23").Value = "  +5.46%  "

# Row 24
Set-TextValue $ws.Range("D24") "11.23"
$ws.Range("E24").Value = "  +1.79%  "

# Row 25
$ws.Range("E25").Value = "  -1.08%  "

# Row 26
Set-TextValue $ws.Range("D26") "2.144.78"
$ws.Range("E26").Value = "  +2.83%  "

# Row 27
Set-TextValue $ws.Range("B27") "Monero"
Set-TextValue $ws.Range("C27") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D27") "160.04"
$ws.Range("E27").Value = "  -1.03%  "

# Row 28
Set-TextValue $ws.Range("B28") "EthereumClassic"
Set-TextValue $ws.Range("C28") "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D28") "21.04"
$ws.Range("E28").Value = "  +2.10%  "

# Row 29
Set-TextValue $ws.Range("D29") "2.471"
$ws.Range("E29").Value = "  +5.90%  "

# Row 30
Set-TextValue $ws.Range("D30") "129.73"
$ws.Range("E30").Value = "  +3.36%  "

# Row 31
Set-TextValue $ws.Range("D31") "1.080"
$ws.Range("E31").Value = "  +3.97%  "

# Row 32
Set-TextValue $ws.Range("D32") "0.1058"
$ws.Range("E32").Value = "  +1.74%  "

# Row 33
Set-TextValue $ws.Range("D33") "6.112"
$ws.Range("E33").Value = "  +6.04%  "

# Row 34
Set-TextValue $ws.Range("D34") "3.664"
$ws.Range("E34").Value = "  +1.65%  "

# Row 35
Set-TextValue $ws.Range("D35") "0.02504"
$ws.Range("E35").Value = "  +2.44%  "

# Row 36
Set-TextValue $ws.Range("D36") "0.06624"
$ws.Range("E36").Value = "  +1.63%  "

# Row 37
Set-TextValue $ws.Range("D37") "0.2214"
$ws.Range("E37").Value = "  +2.74%  "

# Row 38
$ws.Range("E38").Value = "  +5.28%  "

# Row 39
Set-TextValue $ws.Range("D39") "9.043"
$ws.Range("E39").Value = "  +2.41%  "

# Row 40
Set-TextValue $ws.Range("D40") "5.206"
$ws.Range("E40").Value = "  +3.28%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.6530"
$ws.Range("E41").Value = "  +1.97%  "

# Row 42
Set-TextValue $ws.Range("D42") "1.240"
$ws.Range("E42").Value = "  -0.88%  "

# Row 43
Set-TextValue $ws.Range("D43") "11.43"
$ws.Range("E43").Value = "  +3.29%  "

# Row 44
Set-TextValue $ws.Range("D44") "0.6149"
$ws.Range("E44").Value = "  +2.34%  "

# Row 45
Set-TextValue $ws.Range("D45") "13.15"
$ws.Range("E45").Value = "  +1.28%  "

# Row 46
Set-TextValue $ws.Range("D46") "3.723"
$ws.Range("E46").Value = "  +1.01%  "

# Row 47
Set-TextValue $ws.Range("D47") "2.061"
$ws.Range("E47").Value = "  +3.08%  "

# Row 48
Set-TextValue $ws.Range("D48") "1.242"
$ws.Range("E48").Value = "  +2.56%  "

# Row 49
Set-TextValue $ws.Range("D49") "125.58"
$ws.Range("E49").Value = "  +3.23%  "

# Row 50
Set-TextValue $ws.Range("D50") "1.151"
$ws.Range("E50").Value = "  -2.19%  "

# Row 51
Set-TextValue $ws.Range("D51") "79.47"
$ws.Range("E51").Value = "  +4.29%  "
